$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column values are stored as text (inline strings) in the source sheet.
# Prefix with an apostrophe to force text entry (avoids Excel auto-numeric
# conversion of values like "1.016"), then reset the style so no extra
# quote-prefix / number-format style is left behind on the cell.

$ws.Range("D2").Value = "'27.168.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "'1.853.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").Value = "'1.015"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'310.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.4785"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("D8").Value = "'0.3694"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "'0.07266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "'0.9333"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'19.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").Value = "'0.07800"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "'1.803.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "'5.403"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "'6.498"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "'89.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "'0.000008706"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "'27.183.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "'14.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("D22").Value = "'5.071"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'1.950"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").Value = "'153.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "'18.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").Value = "'1.992"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "'114.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "'4.943"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").Value = "'0.08889"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "'3.302"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").Value = "'1.187"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "'4.532"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").Value = "'0.7402"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").Value = "'2.690"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "'1.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("D37").Value = "'0.01986"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").Value = "'0.05279"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "'2.979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "'0.5293"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").Value = "'7.056"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "'0.1527"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'8.318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "'10.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "'0.4755"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "'1.016"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "'102.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "'1.623"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "'65.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").Value = "'0.06068"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "'0.8944"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.69%  "
